$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed scrape numbers (Number Students / Rating Count / Rating Value)
$ws.Range("F2").Value = 153524
$ws.Range("G2").Value = 6241
$ws.Range("F3").Value = 150364
$ws.Range("G3").Value = 31715
$ws.Range("F4").Value = 96548
$ws.Range("G4").Value = 19262
$ws.Range("F5").Value = 58850
$ws.Range("F6").Value = 47517
$ws.Range("G6").Value = 2329
$ws.Range("F7").Value = 39166
$ws.Range("G7").Value = 6514
$ws.Range("F8").Value = 39166
$ws.Range("G8").Value = 6514
$ws.Range("F9").Value = 38760
$ws.Range("G9").Value = 6458
$ws.Range("F11").Value = 32697
$ws.Range("G11").Value = 6955
$ws.Range("F12").Value = 31466
$ws.Range("G12").Value = 5695
$ws.Range("F13").Value = 27316
$ws.Range("F14").Value = 27088
$ws.Range("G14").Value = 4561
$ws.Range("F15").Value = 26009
$ws.Range("G15").Value = 5688
$ws.Range("F16").Value = 23935
$ws.Range("G16").Value = 3835
$ws.Range("F17").Value = 21301
$ws.Range("F18").Value = 20636
$ws.Range("G18").Value = 5291
$ws.Range("H18").Value = 4.5
$ws.Range("F19").Value = 19311
$ws.Range("G19").Value = 3985
$ws.Range("F20").Value = 18490
$ws.Range("F21").Value = 18007
$ws.Range("G21").Value = 4458
$ws.Range("F22").Value = 16241
$ws.Range("G22").Value = 1995
$ws.Range("F23").Value = 14736
$ws.Range("G23").Value = 4766
$ws.Range("F24").Value = 14399
$ws.Range("G24").Value = 3277
$ws.Range("F25").Value = 13999
$ws.Range("F26").Value = 13303
$ws.Range("G26").Value = 2183
$ws.Range("F27").Value = 12916
$ws.Range("G27").Value = 4996
$ws.Range("F30").Value = 9928
$ws.Range("G30").Value = 3447
$ws.Range("F32").Value = 8699
$ws.Range("G32").Value = 2080
$ws.Range("F33").Value = 8253
$ws.Range("G33").Value = 1952
$ws.Range("F34").Value = 5551
$ws.Range("F35").Value = 5157
$ws.Range("F39").Value = 3019
$ws.Range("F40").Value = 2984
$ws.Range("G41").Value = 513
$ws.Range("F43").Value = 2798
$ws.Range("F44").Value = 2639
$ws.Range("H45").Value = 4.6
$ws.Range("F46").Value = 942
$ws.Range("G46").Value = 131
$ws.Range("H46").Value = 4.7
